# actualizacion resumen a mayo 2022
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Apply "0.0" number format to the recently recalculated Inv_Recuperada
#     (column Q) cells for March and April 2022 (rows 51-52), matching the
#     new custom numFmt 164 used for the newly appended rows below.
$ws.Range("Q51").NumberFormat = "0.0"
$ws.Range("Q52").NumberFormat = "0.0"

# --- Row 53: 2022-04 ---------------------------------------------------
$ws.Cells.Item(53, 1).Value2 = 2022
$ws.Cells.Item(53, 2).Value2 = 4
$ws.Cells.Item(53, 3).Value2 = 61993
$ws.Cells.Item(53, 4).Value2 = 102908.38
$ws.Cells.Item(53, 5).Value2 = 92.32
$ws.Cells.Item(53, 6).Value2 = 103000.7
$ws.Cells.Item(53, 7).Value2 = 42663.24
$ws.Cells.Item(53, 8).Value2 = 39531.170000000006
$ws.Cells.Item(53, 9).Value2 = 6917.77
$ws.Cells.Item(53, 10).Value2 = 5847.8
$ws.Cells.Item(53, 11).Value2 = 6558
$ws.Cells.Item(53, 12).Value2 = 101517.98000000001
$ws.Cells.Item(53, 13).Value2 = 1482.7199999999866
$ws.Cells.Item(53, 14).Value2 = 213.03448275861876
$ws.Cells.Item(53, 15).Value2 = 0.023917539076992348
$ws.Cells.Item(53, 16).Value2 = 500000
$ws.Cells.Item(53, 17).Value2 = 270974.65849761362
$ws.Cells.Item(53, 18).Value2 = 0.54194931699522719

$ws.Range("D53:F53").NumberFormat = "0.00"
$ws.Range("J53:K53").NumberFormat = "0.00"
$ws.Range("N53").NumberFormat = "0.00"
$ws.Range("Q53").NumberFormat = "0.0"

# --- Row 54: 2022-05 ---------------------------------------------------
$ws.Cells.Item(54, 1).Value2 = 2022
$ws.Cells.Item(54, 2).Value2 = 5
$ws.Cells.Item(54, 3).Value2 = 61349
$ws.Cells.Item(54, 4).Value2 = 101839.34
$ws.Cells.Item(54, 5).Value2 = 84.99
$ws.Cells.Item(54, 6).Value2 = 101924.33
$ws.Cells.Item(54, 7).Value2 = 38879.949999999997
$ws.Cells.Item(54, 8).Value2 = 40261.050000000003
$ws.Cells.Item(54, 9).Value2 = 7010.62
$ws.Cells.Item(54, 10).Value2 = 3484
$ws.Cells.Item(54, 11).Value2 = 6796
$ws.Cells.Item(54, 12).Value2 = 96431.62
$ws.Cells.Item(54, 13).Value2 = 5492.71
$ws.Cells.Item(54, 14).Value2 = 789.18
$ws.Cells.Item(54, 15).Value2 = 0.089532184713687374
$ws.Cells.Item(54, 16).Value2 = 500000
$ws.Cells.Item(54, 17).Value2 = 271753.78000000003
$ws.Cells.Item(54, 18).Value2 = 0.54350756699522729

$ws.Range("D54:N54").NumberFormat = "0.00"
$ws.Range("O54:P54").NumberFormat = "0.00"
$ws.Range("Q54").NumberFormat = "0.0"

# --- View state: scroll/select near the newly entered data -------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("U40").Select()
